$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.086.46"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.93"
$ws.Range("E3").Value = "  +0.22%  "

$ws.Range("E4").Value = "  -0.40%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.96"
$ws.Range("E5").Value = "  -0.39%  "

$ws.Range("E6").Value = "  -0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4575"
$ws.Range("E7").Value = "  +6.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3730"
$ws.Range("E8").Value = "  +1.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07316"
$ws.Range("E9").Value = "  +1.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8618"
$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.97"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.46"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.702"
$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.10"
$ws.Range("E14").Value = "  +4.37%  "

$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.364"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07101"
$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008848"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.03"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.117.79"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.200"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  +0.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.000"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.86"
$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.220"
$ws.Range("E26").Value = "  +4.40%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.279"
$ws.Range("E28").Value = "  +0.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.41"
$ws.Range("E29").Value = "  +0.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08890"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.199"
$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7608"
$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.972"
$ws.Range("E33").Value = "  +6.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.478"
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.104"
$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01972"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05288"
$ws.Range("E38").Value = "  +0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5371"
$ws.Range("E39").Value = "  +6.73%  "

$ws.Range("E40").Value = "  +1.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.884"
$ws.Range("E41").Value = "  -1.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1716"
$ws.Range("E42").Value = "  +2.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5239"
$ws.Range("E43").Value = "  +11.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.611"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.68"
$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.984"
$ws.Range("E46").Value = "  +9.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.84"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.679"
$ws.Range("E48").Value = "  +0.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06415"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9231"
$ws.Range("E51").Value = "  +1.02%  "
